$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.058.77"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.303.19"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.73"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.33"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.116"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.09"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +10.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.78"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "2.656.19"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "2.293.39"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.813"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "42.957.32"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0901"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.97"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.98"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  +4.05%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.49"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.72"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.13"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +7.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.96"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.88"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0698"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.40"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "1.981.55"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.89"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.25%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.40"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.524.31"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.59"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.74%  "
